# Weekly data refresh: insert two new rows of "Coliflor" price data
# (Vega Central Mapocho de Santiago) at the top of the date-ordered
# block, pushing the existing rows 404-422 down to 406-424.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 404-405; everything currently at row 404
# onward (through 422) shifts down to 406-424.
$ws.Range("A404:A405").EntireRow.Insert()

# New row 404: Primera quality, week of 2021-11-09
$ws.Range("A404").Value = 9
$ws.Range("B404").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C404").Value = "Metropolitana"
$ws.Range("D404").Value = 44509
$ws.Range("E404").Value = 13
$ws.Range("F404").Value = 100112008
$ws.Range("G404").Value = "Coliflor"
$ws.Range("H404").Value = "Sin especificar"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 4300
$ws.Range("K404").Value = 600
$ws.Range("L404").Value = 700
$ws.Range("M404").Value = 650
$ws.Range("N404").Value = "$/unidad"
$ws.Range("O404").Value = "Región Metropolitana"
$ws.Range("P404").Value = 650
$ws.Range("Q404").Value = 1
$ws.Range("R404").Value = "Hortaliza"

# New row 405: Segunda quality, week of 2021-11-09
$ws.Range("A405").Value = 9
$ws.Range("B405").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C405").Value = "Metropolitana"
$ws.Range("D405").Value = 44509
$ws.Range("E405").Value = 13
$ws.Range("F405").Value = 100112008
$ws.Range("G405").Value = "Coliflor"
$ws.Range("H405").Value = "Sin especificar"
$ws.Range("I405").Value = "Segunda"
$ws.Range("J405").Value = 1600
$ws.Range("K405").Value = 500
$ws.Range("L405").Value = 500
$ws.Range("M405").Value = 500
$ws.Range("N405").Value = "$/unidad"
$ws.Range("O405").Value = "Región Metropolitana"
$ws.Range("P405").Value = 500
$ws.Range("Q405").Value = 1
$ws.Range("R405").Value = "Hortaliza"
